# Fix bugs, update isolate column comment if the isolate is not a clinical isolate.
# Capitalize the first letter of each specimen-type segment in column N ("Specimens"),
# leaving "NA" counts untouched (they are already upper-case).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 14).End(-4162).Row  # xlUp, column N = 14
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 14)  # column N
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") { continue }

    $parts = $val.ToString().Split(",")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $piece = $parts[$i]
        $trimmed = $piece.TrimStart(" ")
        $leadLen = $piece.Length - $trimmed.Length
        $lead = $piece.Substring(0, $leadLen)

        if ($trimmed -eq "NA" -or $trimmed.StartsWith("NA ") -or $trimmed.StartsWith("NA(")) {
            # already "NA" style token - leave as-is
            $parts[$i] = $lead + $trimmed
        }
        elseif ($trimmed.Length -gt 0) {
            $capitalized = $trimmed.Substring(0,1).ToUpper() + $trimmed.Substring(1)
            $parts[$i] = $lead + $capitalized
        }
        else {
            $parts[$i] = $piece
        }
    }

    $newVal = [string]::Join(",", $parts)
    if (-not $newVal.Equals($val)) {
        $cell.Value2 = $newVal
    }
}
